$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 108

$ws.Range("E6").Value = 48

$ws.Range("E7").Value = 41
$ws.Range("F7").Value = 25
$ws.Range("H7").Value = 29

$ws.Range("E10").Value = 652
$ws.Range("F10").Value = 349
$ws.Range("H10").Value = 445

$ws.Range("E11").Value = 426
$ws.Range("F11").Value = 232
$ws.Range("H11").Value = 296

$ws.Range("E12").Value = 655
$ws.Range("F12").Value = 381
$ws.Range("H12").Value = 467

$ws.Range("E13").Value = 155
$ws.Range("F13").Value = 85
$ws.Range("H13").Value = 119

$ws.Range("E14").Value = 135

$ws.Range("E15").Value = 190
$ws.Range("F15").Value = 87
$ws.Range("H15").Value = 137

$ws.Range("E16").Value = 226

$ws.Range("E17").Value = 120

$ws.Range("E22").Value = 185

$ws.Range("E23").Value = 217
$ws.Range("F23").Value = 109
$ws.Range("H23").Value = 160

$ws.Range("E25").Value = 308
$ws.Range("F25").Value = 167
$ws.Range("H25").Value = 227

$ws.Range("E26").Value = 180
$ws.Range("F26").Value = 109
$ws.Range("H26").Value = 134

$ws.Range("E27").Value = 367

$ws.Range("F31").Value = 35
$ws.Range("H31").Value = 63

$ws.Range("E32").Value = 203
$ws.Range("F32").Value = 128
$ws.Range("H32").Value = 166

$ws.Range("E34").Value = 242

$ws.Range("E39").Value = 192
$ws.Range("F39").Value = 100
$ws.Range("H39").Value = 151

$ws.Range("E41").Value = 423

$ws.Range("E43").Value = 137
$ws.Range("F43").Value = 74
$ws.Range("H43").Value = 101

$ws.Range("E44").Value = 348
$ws.Range("F44").Value = 180
$ws.Range("H44").Value = 248

$ws.Range("E45").Value = 171
$ws.Range("F45").Value = 90
$ws.Range("H45").Value = 129

$ws.Range("E46").Value = 370
$ws.Range("F46").Value = 203
$ws.Range("H46").Value = 266

$ws.Range("E47").Value = 516

$ws.Range("E48").Value = 251
$ws.Range("F48").Value = 115
$ws.Range("H48").Value = 159

$ws.Range("E50").Value = 267

$ws.Range("E51").Value = 258
